$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 0.507
$ws.Cells.Item(2, 7).Value = 0.007
$ws.Cells.Item(2, 8).Value = 0.227
$ws.Cells.Item(2, 9).Value = 0.127
$ws.Cells.Item(2, 10).Value = 0.13
$ws.Cells.Item(2, 12).Value = 0.00007
$ws.Cells.Item(3, 6).Value = 0.517
$ws.Cells.Item(3, 7).Value = 0.007
$ws.Cells.Item(3, 8).Value = 0.238
$ws.Cells.Item(3, 9).Value = 0.132
$ws.Cells.Item(3, 10).Value = 0.125
$ws.Cells.Item(3, 12).Value = 0.00007
$ws.Cells.Item(4, 6).Value = 0.515
$ws.Cells.Item(4, 7).Value = 0.007
$ws.Cells.Item(4, 8).Value = 0.235
$ws.Cells.Item(4, 9).Value = 0.13
$ws.Cells.Item(4, 10).Value = 0.13
$ws.Cells.Item(4, 12).Value = 0.00007
$ws.Cells.Item(5, 6).Value = 0.518
$ws.Cells.Item(5, 7).Value = 0.008
$ws.Cells.Item(5, 8).Value = 0.235
$ws.Cells.Item(5, 9).Value = 0.129
$ws.Cells.Item(5, 10).Value = 0.133
$ws.Cells.Item(5, 12).Value = 0.00006
$ws.Cells.Item(6, 6).Value = 0.515
$ws.Cells.Item(6, 7).Value = 0.007
$ws.Cells.Item(6, 8).Value = 0.237
$ws.Cells.Item(6, 9).Value = 0.129
$ws.Cells.Item(6, 10).Value = 0.13
$ws.Cells.Item(6, 12).Value = 0.00006
$ws.Cells.Item(7, 6).Value = 0.558
$ws.Cells.Item(7, 7).Value = 0.007
$ws.Cells.Item(7, 8).Value = 0.24
$ws.Cells.Item(7, 9).Value = 0.17
$ws.Cells.Item(7, 10).Value = 0.127
$ws.Cells.Item(7, 12).Value = 0.00007
$ws.Cells.Item(8, 6).Value = 0.509
$ws.Cells.Item(8, 7).Value = 0.008
$ws.Cells.Item(8, 8).Value = 0.236
$ws.Cells.Item(8, 9).Value = 0.126
$ws.Cells.Item(8, 10).Value = 0.126
$ws.Cells.Item(8, 12).Value = 0.00007
$ws.Cells.Item(9, 6).Value = 0.518
$ws.Cells.Item(9, 7).Value = 0.007
$ws.Cells.Item(9, 8).Value = 0.237
$ws.Cells.Item(9, 9).Value = 0.128
$ws.Cells.Item(9, 10).Value = 0.133
$ws.Cells.Item(9, 12).Value = 0.00006
$ws.Cells.Item(10, 6).Value = 0.498
$ws.Cells.Item(10, 7).Value = 0.007
$ws.Cells.Item(10, 8).Value = 0.228
$ws.Cells.Item(10, 9).Value = 0.125
$ws.Cells.Item(10, 10).Value = 0.125
$ws.Cells.Item(10, 12).Value = 0.00006
$ws.Cells.Item(11, 6).Value = 0.504
$ws.Cells.Item(11, 7).Value = 0.007
$ws.Cells.Item(11, 8).Value = 0.234
$ws.Cells.Item(11, 9).Value = 0.125
$ws.Cells.Item(11, 10).Value = 0.124
$ws.Cells.Item(11, 12).Value = 0.00007
$ws.Cells.Item(12, 6).Value = 0.516
$ws.Cells.Item(12, 7).Value = 0.008
$ws.Cells.Item(12, 8).Value = 0.234
$ws.Cells.Item(12, 9).Value = 0.131
$ws.Cells.Item(12, 10).Value = 0.131
$ws.Cells.Item(12, 12).Value = 0.00007
$ws.Cells.Item(13, 6).Value = 0.51
$ws.Cells.Item(13, 7).Value = 0.007
$ws.Cells.Item(13, 8).Value = 0.234
$ws.Cells.Item(13, 9).Value = 0.128
$ws.Cells.Item(13, 10).Value = 0.128
$ws.Cells.Item(13, 12).Value = 0.00006
$ws.Cells.Item(14, 6).Value = 0.517
$ws.Cells.Item(14, 7).Value = 0.007
$ws.Cells.Item(14, 8).Value = 0.239
$ws.Cells.Item(14, 9).Value = 0.129
$ws.Cells.Item(14, 10).Value = 0.129
$ws.Cells.Item(14, 12).Value = 0.00007
$ws.Cells.Item(15, 6).Value = 0.507
$ws.Cells.Item(15, 7).Value = 0.007
$ws.Cells.Item(15, 8).Value = 0.232
$ws.Cells.Item(15, 9).Value = 0.124
$ws.Cells.Item(15, 10).Value = 0.13
$ws.Cells.Item(15, 12).Value = 0.00007
$ws.Cells.Item(16, 6).Value = 0.513
$ws.Cells.Item(16, 7).Value = 0.008
$ws.Cells.Item(16, 8).Value = 0.233
$ws.Cells.Item(16, 9).Value = 0.125
$ws.Cells.Item(16, 10).Value = 0.134
$ws.Cells.Item(16, 12).Value = 0.00007

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 1.053
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0.002
$ws.Cells.Item(2, 9).Value = 0.352
$ws.Cells.Item(2, 10).Value = 0.355
$ws.Cells.Item(2, 12).Value = 0.00006
$ws.Cells.Item(3, 6).Value = 1.057
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0.002
$ws.Cells.Item(3, 9).Value = 0.347
$ws.Cells.Item(3, 10).Value = 0.37
$ws.Cells.Item(3, 12).Value = 0.00007
$ws.Cells.Item(4, 6).Value = 1.05
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0.002
$ws.Cells.Item(4, 9).Value = 0.354
$ws.Cells.Item(4, 10).Value = 0.358
$ws.Cells.Item(4, 12).Value = 0.00006
$ws.Cells.Item(5, 6).Value = 1.045
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0.002
$ws.Cells.Item(5, 9).Value = 0.351
$ws.Cells.Item(5, 10).Value = 0.362
$ws.Cells.Item(5, 12).Value = 0.00007
$ws.Cells.Item(6, 6).Value = 1.043
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0.002
$ws.Cells.Item(6, 9).Value = 0.355
$ws.Cells.Item(6, 10).Value = 0.355
$ws.Cells.Item(6, 12).Value = 0.00006
$ws.Cells.Item(7, 6).Value = 1.05
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0.002
$ws.Cells.Item(7, 9).Value = 0.353
$ws.Cells.Item(7, 10).Value = 0.353
$ws.Cells.Item(7, 12).Value = 0.00007
$ws.Cells.Item(8, 6).Value = 1.034
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0.002
$ws.Cells.Item(8, 9).Value = 0.349
$ws.Cells.Item(8, 10).Value = 0.347
$ws.Cells.Item(8, 12).Value = 0.00007
$ws.Cells.Item(9, 6).Value = 1.04
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0.002
$ws.Cells.Item(9, 9).Value = 0.346
$ws.Cells.Item(9, 10).Value = 0.354
$ws.Cells.Item(9, 12).Value = 0.00007
$ws.Cells.Item(10, 6).Value = 1.064
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0.002
$ws.Cells.Item(10, 9).Value = 0.364
$ws.Cells.Item(10, 10).Value = 0.368
$ws.Cells.Item(10, 12).Value = 0.00007
$ws.Cells.Item(11, 6).Value = 1.051
$ws.Cells.Item(11, 7).Value = 0.001
$ws.Cells.Item(11, 8).Value = 0.002
$ws.Cells.Item(11, 9).Value = 0.355
$ws.Cells.Item(11, 10).Value = 0.356
$ws.Cells.Item(11, 12).Value = 0.00007
$ws.Cells.Item(12, 6).Value = 1.077
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0.002
$ws.Cells.Item(12, 9).Value = 0.376
$ws.Cells.Item(12, 10).Value = 0.354
$ws.Cells.Item(12, 12).Value = 0.00007
$ws.Cells.Item(13, 6).Value = 1.051
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0.002
$ws.Cells.Item(13, 9).Value = 0.351
$ws.Cells.Item(13, 10).Value = 0.355
$ws.Cells.Item(13, 12).Value = 0.00007
$ws.Cells.Item(14, 6).Value = 1.038
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0.002
$ws.Cells.Item(14, 9).Value = 0.345
$ws.Cells.Item(14, 10).Value = 0.352
$ws.Cells.Item(14, 12).Value = 0.00007
$ws.Cells.Item(15, 6).Value = 1.043
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0.002
$ws.Cells.Item(15, 9).Value = 0.357
$ws.Cells.Item(15, 10).Value = 0.352
$ws.Cells.Item(15, 12).Value = 0.00006
$ws.Cells.Item(16, 6).Value = 1.036
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0.002
$ws.Cells.Item(16, 9).Value = 0.349
$ws.Cells.Item(16, 10).Value = 0.358
$ws.Cells.Item(16, 12).Value = 0.00006

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 0.507
$ws.Cells.Item(2, 7).Value = 0.007
$ws.Cells.Item(2, 8).Value = 0.227
$ws.Cells.Item(2, 9).Value = 0.127
$ws.Cells.Item(2, 10).Value = 0.13
$ws.Cells.Item(2, 12).Value = 0.00007
$ws.Cells.Item(3, 6).Value = 0.517
$ws.Cells.Item(3, 7).Value = 0.007
$ws.Cells.Item(3, 8).Value = 0.238
$ws.Cells.Item(3, 9).Value = 0.132
$ws.Cells.Item(3, 10).Value = 0.125
$ws.Cells.Item(3, 12).Value = 0.00007
$ws.Cells.Item(4, 6).Value = 0.515
$ws.Cells.Item(4, 7).Value = 0.007
$ws.Cells.Item(4, 8).Value = 0.235
$ws.Cells.Item(4, 9).Value = 0.13
$ws.Cells.Item(4, 10).Value = 0.13
$ws.Cells.Item(4, 12).Value = 0.00007
$ws.Cells.Item(5, 6).Value = 0.518
$ws.Cells.Item(5, 7).Value = 0.008
$ws.Cells.Item(5, 8).Value = 0.235
$ws.Cells.Item(5, 9).Value = 0.129
$ws.Cells.Item(5, 10).Value = 0.133
$ws.Cells.Item(5, 12).Value = 0.00006
$ws.Cells.Item(6, 6).Value = 0.515
$ws.Cells.Item(6, 7).Value = 0.007
$ws.Cells.Item(6, 8).Value = 0.237
$ws.Cells.Item(6, 9).Value = 0.129
$ws.Cells.Item(6, 10).Value = 0.13
$ws.Cells.Item(6, 12).Value = 0.00006
$ws.Cells.Item(7, 6).Value = 0.558
$ws.Cells.Item(7, 7).Value = 0.007
$ws.Cells.Item(7, 8).Value = 0.24
$ws.Cells.Item(7, 9).Value = 0.17
$ws.Cells.Item(7, 10).Value = 0.127
$ws.Cells.Item(7, 12).Value = 0.00007
$ws.Cells.Item(8, 6).Value = 0.509
$ws.Cells.Item(8, 7).Value = 0.008
$ws.Cells.Item(8, 8).Value = 0.236
$ws.Cells.Item(8, 9).Value = 0.126
$ws.Cells.Item(8, 10).Value = 0.126
$ws.Cells.Item(8, 12).Value = 0.00007
$ws.Cells.Item(9, 6).Value = 0.518
$ws.Cells.Item(9, 7).Value = 0.007
$ws.Cells.Item(9, 8).Value = 0.237
$ws.Cells.Item(9, 9).Value = 0.128
$ws.Cells.Item(9, 10).Value = 0.133
$ws.Cells.Item(9, 12).Value = 0.00006
$ws.Cells.Item(10, 6).Value = 0.498
$ws.Cells.Item(10, 7).Value = 0.007
$ws.Cells.Item(10, 8).Value = 0.228
$ws.Cells.Item(10, 9).Value = 0.125
$ws.Cells.Item(10, 10).Value = 0.125
$ws.Cells.Item(10, 12).Value = 0.00006
$ws.Cells.Item(11, 6).Value = 0.504
$ws.Cells.Item(11, 7).Value = 0.007
$ws.Cells.Item(11, 8).Value = 0.234
$ws.Cells.Item(11, 9).Value = 0.125
$ws.Cells.Item(11, 10).Value = 0.124
$ws.Cells.Item(11, 12).Value = 0.00007
$ws.Cells.Item(12, 6).Value = 0.516
$ws.Cells.Item(12, 7).Value = 0.008
$ws.Cells.Item(12, 8).Value = 0.234
$ws.Cells.Item(12, 9).Value = 0.131
$ws.Cells.Item(12, 10).Value = 0.131
$ws.Cells.Item(12, 12).Value = 0.00007
$ws.Cells.Item(13, 6).Value = 0.51
$ws.Cells.Item(13, 7).Value = 0.007
$ws.Cells.Item(13, 8).Value = 0.234
$ws.Cells.Item(13, 9).Value = 0.128
$ws.Cells.Item(13, 10).Value = 0.128
$ws.Cells.Item(13, 12).Value = 0.00006
$ws.Cells.Item(14, 6).Value = 0.517
$ws.Cells.Item(14, 7).Value = 0.007
$ws.Cells.Item(14, 8).Value = 0.239
$ws.Cells.Item(14, 9).Value = 0.129
$ws.Cells.Item(14, 10).Value = 0.129
$ws.Cells.Item(14, 12).Value = 0.00007
$ws.Cells.Item(15, 6).Value = 0.507
$ws.Cells.Item(15, 7).Value = 0.007
$ws.Cells.Item(15, 8).Value = 0.232
$ws.Cells.Item(15, 9).Value = 0.124
$ws.Cells.Item(15, 10).Value = 0.13
$ws.Cells.Item(15, 12).Value = 0.00007
$ws.Cells.Item(16, 6).Value = 0.513
$ws.Cells.Item(16, 7).Value = 0.008
$ws.Cells.Item(16, 8).Value = 0.233
$ws.Cells.Item(16, 9).Value = 0.125
$ws.Cells.Item(16, 10).Value = 0.134
$ws.Cells.Item(16, 12).Value = 0.00007
$ws.Cells.Item(17, 6).Value = 1.053
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0.002
$ws.Cells.Item(17, 9).Value = 0.352
$ws.Cells.Item(17, 10).Value = 0.355
$ws.Cells.Item(17, 12).Value = 0.00006
$ws.Cells.Item(18, 6).Value = 1.057
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0.002
$ws.Cells.Item(18, 9).Value = 0.347
$ws.Cells.Item(18, 10).Value = 0.37
$ws.Cells.Item(18, 12).Value = 0.00007
$ws.Cells.Item(19, 6).Value = 1.05
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0.002
$ws.Cells.Item(19, 9).Value = 0.354
$ws.Cells.Item(19, 10).Value = 0.358
$ws.Cells.Item(19, 12).Value = 0.00006
$ws.Cells.Item(20, 6).Value = 1.045
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0.002
$ws.Cells.Item(20, 9).Value = 0.351
$ws.Cells.Item(20, 10).Value = 0.362
$ws.Cells.Item(20, 12).Value = 0.00007
$ws.Cells.Item(21, 6).Value = 1.043
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0.002
$ws.Cells.Item(21, 9).Value = 0.355
$ws.Cells.Item(21, 10).Value = 0.355
$ws.Cells.Item(21, 12).Value = 0.00006
$ws.Cells.Item(22, 6).Value = 1.05
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0.002
$ws.Cells.Item(22, 9).Value = 0.353
$ws.Cells.Item(22, 10).Value = 0.353
$ws.Cells.Item(22, 12).Value = 0.00007
$ws.Cells.Item(23, 6).Value = 1.034
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0.002
$ws.Cells.Item(23, 9).Value = 0.349
$ws.Cells.Item(23, 10).Value = 0.347
$ws.Cells.Item(23, 12).Value = 0.00007
$ws.Cells.Item(24, 6).Value = 1.04
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 0.002
$ws.Cells.Item(24, 9).Value = 0.346
$ws.Cells.Item(24, 10).Value = 0.354
$ws.Cells.Item(24, 12).Value = 0.00007
$ws.Cells.Item(25, 6).Value = 1.064
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0.002
$ws.Cells.Item(25, 9).Value = 0.364
$ws.Cells.Item(25, 10).Value = 0.368
$ws.Cells.Item(25, 12).Value = 0.00007
$ws.Cells.Item(26, 6).Value = 1.051
$ws.Cells.Item(26, 7).Value = 0.001
$ws.Cells.Item(26, 8).Value = 0.002
$ws.Cells.Item(26, 9).Value = 0.355
$ws.Cells.Item(26, 10).Value = 0.356
$ws.Cells.Item(26, 12).Value = 0.00007
$ws.Cells.Item(27, 6).Value = 1.077
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 0.002
$ws.Cells.Item(27, 9).Value = 0.376
$ws.Cells.Item(27, 10).Value = 0.354
$ws.Cells.Item(27, 12).Value = 0.00007
$ws.Cells.Item(28, 6).Value = 1.051
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0.002
$ws.Cells.Item(28, 9).Value = 0.351
$ws.Cells.Item(28, 10).Value = 0.355
$ws.Cells.Item(28, 12).Value = 0.00007
$ws.Cells.Item(29, 6).Value = 1.038
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 0.002
$ws.Cells.Item(29, 9).Value = 0.345
$ws.Cells.Item(29, 10).Value = 0.352
$ws.Cells.Item(29, 12).Value = 0.00007
$ws.Cells.Item(30, 6).Value = 1.043
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0.002
$ws.Cells.Item(30, 9).Value = 0.357
$ws.Cells.Item(30, 10).Value = 0.352
$ws.Cells.Item(30, 12).Value = 0.00006
$ws.Cells.Item(31, 6).Value = 1.036
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0.002
$ws.Cells.Item(31, 9).Value = 0.349
$ws.Cells.Item(31, 10).Value = 0.358
$ws.Cells.Item(31, 12).Value = 0.00006

Write-Output "done"